$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15, column A: tiny floating point precision correction to the existing date-time value
$ws.Cells.Item(15, 1).Value = 44328.77718483681

# New row 16 data
$ws.Cells.Item(16, 1).Value = 44329.77915936564
$ws.Cells.Item(16, 2).Value = 74804
$ws.Cells.Item(16, 3).Value = 62922
$ws.Cells.Item(16, 4).Value = 3177
$ws.Cells.Item(16, 5).Value = 2106
$ws.Cells.Item(16, 6).Value = 1491
$ws.Cells.Item(16, 7).Value = 19460
$ws.Cells.Item(16, 8).Value = 1292
$ws.Cells.Item(16, 9).Value = 868
$ws.Cells.Item(16, 10).Value = 202
